$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.534.98'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '3.457.36'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.03'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.81'
$ws.Range("E6").Value = '  +4.74%  '
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '3.457.97'
$ws.Range("E8").Value = '  +1.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.586'
$ws.Range("E9").Value = '  +9.83%  '
$ws.Range("E10").Value = '  -1.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.126'
$ws.Range("E11").Value = '  +4.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  +2.23%  '
$ws.Range("D13").Value = '4.057.78'
$ws.Range("E13").Value = '  +2.15%  '
$ws.Range("E14").Value = '  -2.80%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.09'
$ws.Range("E15").Value = '  +8.08%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000194'
$ws.Range("E16").Value = '  +5.99%  '
$ws.Range("D17").Value = '64.555.81'
$ws.Range("E17").Value = '  +1.87%  '
$ws.Range("D18").Value = '3.455.62'
$ws.Range("E18").Value = '  +2.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.41'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.49'
$ws.Range("E20").Value = '  +3.58%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.57'
$ws.Range("E21").Value = '  +1.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.22'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.549'
$ws.Range("E23").Value = '  +2.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.19'
$ws.Range("E24").Value = '  +3.52%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +20.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.19'
$ws.Range("E30").Value = '  +10.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  +9.31%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.57'
$ws.Range("E33").Value = '  +1.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.69'
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.12'
$ws.Range("E36").Value = '  +5.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.51'
$ws.Range("E37").Value = '  +2.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.64'
$ws.Range("E38").Value = '  +1.92%  '
$ws.Range("E39").Value = '  +1.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0776'
$ws.Range("E40").Value = '  +2.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.60'
$ws.Range("E41").Value = '  +0.60%  '
$ws.Range("D42").Value = '2.927.43'
$ws.Range("E42").Value = '  +1.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.56'
$ws.Range("E43").Value = '  +6.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0319'
$ws.Range("E44").Value = '  -0.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.73'
$ws.Range("E45").Value = '  +3.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.773'
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.27'
$ws.Range("E47").Value = '  +9.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("E48").Value = '  +3.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("E49").Value = '  +15.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.872'
$ws.Range("E50").Value = '  +7.20%  '
$ws.Range("E51").Value = '  +4.56%  '
